$ErrorActionPreference = "Stop"
$d = $word.ActiveDocument

function Replace-UniqueText($searchText, $replaceText) {
    $r = $d.Content.Duplicate
    $ok = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Text not found (unique): $searchText"
    }
    $r.Text = $replaceText
}

function Replace-TextInParagraph($paragraphIndex, $searchText, $replaceText) {
    $p = $d.Paragraphs.Item($paragraphIndex)
    $r = $p.Range.Duplicate
    $ok = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Text not found in paragraph $paragraphIndex : $searchText"
    }
    $r.Text = $replaceText
}

function Replace-BoldRunInParagraph($paragraphIndex, $searchText, $replaceText) {
    $p = $d.Paragraphs.Item($paragraphIndex)
    $searchRange = $p.Range.Duplicate
    while ($true) {
        $r = $searchRange.Duplicate
        $ok = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if (-not $ok) {
            throw "Bold text not found in paragraph $paragraphIndex : $searchText"
        }
        if ($r.Bold -eq -1) {
            $r.Text = $replaceText
            return
        } else {
            $searchRange.Start = $r.End
            $searchRange.End = $p.Range.End
        }
    }
}


# --- Hunk 0 ---
Replace-UniqueText 'Appendix 11: SWIFT Child Safety Risk and Exposure Interview Guide: Adolescents' 'Bylaag 11: SWIFT Kinderveiligheidsrisiko en Blootstelling Onderhoudgids: Adolessente'
# --- Hunk 1 ---
Replace-UniqueText 'Briefing and telephonic assent:' 'Inligtingsessie en telefoniese toestemming:'
# --- Hunk 2 ---
Replace-UniqueText 'Hi there _____, I’m ________. Your X [mom/dad/granny etc.] said that this would be a good time to talk to you. Is it?' 'Hallo daar _____, Ek is ________. Jou X [ma/pa/ouma ens.] het gesê dit is ''n goeie tyd om met jou te praat. Is dit?'
# --- Hunk 3 ---
Replace-TextInParagraph 6 'If no - When would you like to talk to us? We’ll need about 20 minutes. ' 'Indien nee - Wanneer sal jy graag met ons wil praat? Ons het so 20 minute nodig. '
Replace-TextInParagraph 6 'Book when it’s possible ' 'Bespreek wanneer dit moontlik is '
# --- Hunk 4 ---
Replace-UniqueText '[* only begin only once participant has answered yes and the time is as arranged]' '[* begin slegs wanneer die deelnemer "ja" geantwoord het en die tyd gereël is]'
# --- Hunk 5 ---
Replace-UniqueText 'Thank you for making the time to speak to me. You would have spoken to one of our team already, when they asked you a few questions about your life and your relationship with your X (parent/caregiver role). I wanted to know if we would be able to speak to you again about this but without the “Never”, “Sometimes”, “Often” answers. This time it''s just going to be like a conversation. ' 'Dankie dat jy tyd maak om met my te praat. Jy sou alklaar met een van ons spanlede gepraat het, wanneer hulle ''n paar vrae gevra het oor jou lewe en jou verhouding met jou X (ouer/versorger rol). Ek wou weet of ons weer met jou kan gesels daaroor, maar sonder die "Nooit", "Soms", "Dikwels" antwoorde. Die keer gaan net soos ''n gesprek wees. '
# --- Hunk 6 ---
Replace-TextInParagraph 9 'The same as last time, we got permission to speak to you from your X, but even though they have given permission, I want to make sure ' 'Net soos die laaste keer, het ons toestemming gekry om met jou te praat van jou X, maar al het hulle toestemming gegee, wil ek seker maak '
Replace-BoldRunInParagraph 9 'you' ' jy '
Replace-TextInParagraph 9 ' are okay with speaking to me again. Before you decide if you are okay to talk to me, I want to tell you a little bit more about everything like the last time so that you know what you are saying yes to. As I explain things I will keep asking if you understand, please let me know if I was being confusing or used a word or spoke about something you didn’t understand. Can I go ahead and explain?' ' is reg daarmee om weer met my te praat. Voordat jy besluit of jy reg is om met my te praat, wil ek vir jou ''n bietjie meer vertel, soos laas, sodat jy weet waarvoor jy ja sê. Soos ek dinge verduidelik, sal ek gereeld vra of jy verstaan, laat weet my asseblief of ek verwarrend is of ''n woord gebruik of oor iets praat wat jy nie verstaan nie. Kan ek voortgaan en verduidelik?'
# --- Hunk 7 ---
Replace-TextInParagraph 10 'This phone call will only take us about 30 minutes. I’m not sure if you remember but I’m working on a project that offers support to parents and caregivers through WhatsApp. We are looking at how it''s working and how the children have found it. We want to hear what your experience of it was. There are no right or wrong answers here. Just be honest, we just want to hear what ' 'Die oproep sal ons net ongeveer 30 minute neem. Ek is nie seker of jy onthou nie, maar ek werk op ''n projek wat ondersteuning aan ouers en versorgers aanbied deur WhatsApp. Ons kyk na hoe dit werk en hoe die kinders dit ervaar het. Ons wil graag hoor wat jou ervaring daarmee was. Daar is geen regte of verkeerde antwoorde hier nie. Wees net eerlik, ons wil net hoor wat '
Replace-BoldRunInParagraph 10 'you' 'jy'
Replace-TextInParagraph 10 ' think. Your family won’t know what you have answered. We will only share something with them if we are worried about your safety and we will first let you know and speak to you about it before we share it. Do you have any questions about this?' ' dink. Jou familie sal nie weet wat jy geantwoord het nie. Ons sal net iets met hulle deel as ons bekommerd is oor jou veiligheid en ons sal jou eers laat weet en met jou praat voordat ons dit deel. Het jy enige vrae hieroor?'
# --- Hunk 8 ---
Replace-UniqueText 'I am working with other people on this project. I want to share the important things that you tell me but when I share it, I will give you a different name. We want to keep your name and identity secret so we will use another name for you when we share anything. Would you like to choose the name we use for you?   ' 'Ek werk met ander mense op die projek. Ek wil graag die belangrike dinge wat jy my vertel deel, maar wanneer ek dit deel, sal ek jou ''n ander naam gee. Ons wil jou naam en identiteit geheim hou, so ons sal ''n ander naam vir jou gebruik wanneer ons iets deel. Wil jy die naam kies wat ons vir jou sal gebruik?   '
# --- Hunk 9 ---
Replace-UniqueText 'If it''s okay with you I will be recording this phone call, so I can remember your answers later. Only the people working with me on this project will hear what you shared but we won’t share your real name with anyone. Do I have permission to record? ' 'As dit reg is met jou, sal ek hierdie oproep opneem sodat ek jou antwoorde later kan onthou. Net die mense wat saam met my aan hierdie projek werk, sal hoor wat jy sê, maar ons sal jou regte naam met niemand deel nie. Het ek jou toestemming om op te neem? '
# --- Hunk 10 ---
Replace-UniqueText 'Even if you give me permission to have this conversation with you, if you change your mind at any point or if you don’t want to answer something specific then you can just let me know. ' 'Selfs as jy vir my toestemming gegee het om met jou te praat, as jy jou mening op enige punt verander het of as jy nie iets spesifiek wil antwoord nie, kan jy net vir my laat weet. '
# --- Hunk 11 ---
Replace-UniqueText 'Do you have any questions? ' 'Het jy enige vrae? '
# --- Hunk 12 ---
Replace-UniqueText 'Can I have your permission to have this conversation with you?' 'Mag ek jou toestemming kry om hierdie gesprek met jou te hê?'
# --- Hunk 13 ---
Replace-UniqueText 'If yes* - Thank you so much, can we begin?' 'Indien ja* – Baie dankie, kan ons begin?'
# --- Hunk 14 ---
Replace-UniqueText 'If they say no - no problem at all, thanks for listening to me. ' 'As hulle nee sê – geen probleem nie, dankie dat jy na my geluister het. '
# --- Hunk 15 ---
Replace-TextInParagraph 19 'Thanks again for making the time to talk to me. I want to make sure that you are in a place where you feel safe to talk. Are you in a place where you feel like you can talk without being overheard by anyone who you don’t want to hear your answers? ' 'Weereens dankie dat jy die tyd geneem het om met my te gesels. Ek wil seker maak jy is in ''n plek waar jy veilig voel om te praat. Is jy op ''n plek waar jy voel jy kan praat sonder dat iemand jou afluister, veral iemand wie jy nie wil hê moet jou antwoorde hoor nie? '
Replace-TextInParagraph 19 '(wait for answer) ' '(wag vir antwoord) '
# --- Hunk 16 ---
Replace-TextInParagraph 20 'If they aren’t' 'Indien nie'
Replace-TextInParagraph 20 ': Would you like to move?' ': Wil jy graag skuif?'
# --- Hunk 17 ---
Replace-TextInParagraph 21 'If they are or once they have moved: ' 'As hulle is of wanneer hulle geskuif het: '
Replace-TextInParagraph 21 'If someone comes into the room you can just tell me to pause and I will wait until you let me know when it''s okay to continue. ' ' As iemand in die kamer kom kan jy net vir my sê om ''n blaaskans te vat en ek sal wag tot jy vir my laat weet dit is reg om voort te gaan. '
# --- Hunk 18 ---
Replace-UniqueText 'Home life and relationship with User' 'Huis lewe en verhouding met Gebruiker'
# --- Hunk 19 ---
Replace-UniqueText 'Can you tell me a bit about your family?' 'Kan jy my bietjie oor jou familie vertel?'
# --- Hunk 20 ---
Replace-UniqueText 'Who lives at home with you? ' 'Wie bly by die huis saam met jou? '
# --- Hunk 21 ---
Replace-UniqueText 'Who is the main person who looks after you? - ' 'Wie is die hoof persoon wat na jou kyk? - '
# --- Hunk 22 ---
Replace-UniqueText 'Probe - Is this the person who has been working through the program?' 'Ondersoek - Is dit die persoon wat deur die program gewerk het?'
# --- Hunk 23 ---
Replace-UniqueText 'Probe - If not, who is the person who worked through the program to you? What do you call them? (Interviewer to refer to X as this for the rest of the interview, e.g. your dad) ' 'Ondersoek - Indien nie, wie is die persoon wat deur die program met jou gewerk het? Wat noem jy hulle? (Onderhoudvoerder moet soos dit na X verwys vir die res van die onderhoud, bv. jou pa) '
# --- Hunk 24 ---
Replace-UniqueText 'Did X tell you that they were working through a parenting program on Whatsapp? (If not, remind them about what ParentText was about)' 'Het X vir jou vertel dat hulle deur ''n ouderskapprogram werk op Whatsapp? (Indien nie, herinner hulle waaroor ParentText gaan)'
# --- Hunk 25 ---
Replace-UniqueText 'How did they tell you about it? What did they say about it? ' 'Hoe het hulle jou daaroor vertel? Wat het hulle daaroor gesê? '
# --- Hunk 26 ---
Replace-UniqueText 'Did your X ever show you the WhatsApp lessons on their phone? If yes - What parts do you remember them showing you?' 'Het jou X jou ooit die WhatsApp lesse op hulle selfoon gewys? Indien ja - Watter dele kan jy onthou van dit wat hulle vir jou gewys het?'
# --- Hunk 27 ---
Replace-UniqueText 'There were some homework activities that they did with you. Which of these homework activities do you remember? Any others? ' 'Daar was ''n paar huiswerk-aktiwiteite wat hulle saam met jou gedoen het. Watter van die huiswerk-aktiwiteite onthou jy? Enige ander? '
# --- Hunk 28 ---
Replace-UniqueText 'Probe: Do they still do any of these activities?' 'Ondersoek: Doen julle nog enige van hierdie aktiwiteite?'
# --- Hunk 29 ---
Replace-UniqueText 'What is your relationship with X like? ' 'Hoe is jou verhouding met X? '
# --- Hunk 30 ---
Replace-UniqueText 'Was it always like this? Did you notice any changes since they started working through the parenting program?' 'Was dit altyd so? Kon jy enige veranderinge op let vanaf hulle deur die ouderskapprogram begin werk het?'
# --- Hunk 31 ---
Replace-UniqueText 'Probe: What were the changes they noticed? What''s improved? What, if anything, got worse or more difficult?' 'Ondersoek: Wat is die veranderinge wat jy opgemerk het? Wat het verbeter? Wat, indien enige, het vererger of moeiliker geraak?'
# --- Hunk 32 ---
Replace-UniqueText 'Since X did the program, have you spent more special time together? What sorts of things do you do together with your X during this time, which you didn’t do before they worked through the WhatsApp programme? added anything to what you do together that is your special time together? (probing quality time)' 'Vandat X die program begin het, het julle meer spesiale tyd saam spander? Watter tipe dinge doen jy saam met jou X gedurende daardie tyd wat julle nie gedoen het voordat julle die WhatsApp-program voltooi het nie? sluit in enige iets wat julle saam doen wat julle spesiale tyd is? (ondersoek kwaliteit tyd)'
# --- Hunk 33 ---
Replace-TextInParagraph 42 'Since X did the program do ' 'Sedert X die program doen, vind '
Replace-BoldRunInParagraph 42 'you' 'jy'
Replace-TextInParagraph 42 ' find it any easier to talk to them about things that worry you? ' ' dit makliker om met hulle te praat oor dinge wat jou bekommer? '
# --- Hunk 34 ---
Replace-TextInParagraph 43 'Since X did the program does it seem like ' 'Vandat X die program gedoen het, lyk dit of '
Replace-BoldRunInParagraph 43 'they ' 'hulle'
Replace-TextInParagraph 43 'find it easier to talk to you about difficult things?' 'dit makliker vind om met jou te gesels oor moeilike dinge?'
# --- Hunk 35 ---
Replace-UniqueText '4.  Did X ever show you the help menu in the whatsapp programme? ' '4.  Het X ooit vir jou die hulp kieslys gewys op die WhatsApp-program? '
